# Corrección en diccionario de datos
# Applies text/content corrections to the data dictionary sheet, fixes a couple of
# PK/FK flags that were swapped, adjusts a few row heights / column width, and
# updates the view (zoom + selection) to match the saved state after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix description text that was wrong / duplicated across several rows
# ---------------------------------------------------------------------------

# INSTITUCION.nombre: description used to (incorrectly) read "Nombre de la
# categoría" - correct it to refer to the institution.
$ws.Range("E17").Value = "Nombre de la institución"

# INSTITUCION.calificacion: clarify the description.
$ws.Range("E20").Value = "Calificación de la institución dentro de la plataforma, indica su reputación."

# ESTUDIANTE.* columns: make wording consistent ("del estudiante" instead of
# "estudiante").
$ws.Range("E25").Value = "Nombre del estudiante"
$ws.Range("E26").Value = "Apellido del estudiante"
$ws.Range("E27").Value = "Usuario de ingreso del estudiante"
$ws.Range("E29").Value = "nacionalidad del estudiante"

# CURSO.nombre: used to wrongly re-use the "Nombre estudiante" text.
$ws.Range("E35").Value = "Nombre del curso"

# CURSO.id_institucion: PK/FK flags were swapped (it's a FK, not a PK).
$ws.Range("F36").Value = "No"
$ws.Range("G36").Value = "Si"

# ESPECIALIZACION.nombre: used to wrongly re-use the "Nombre estudiante" text.
$ws.Range("E44").Value = "Nombre de la especialización"

# ESPECIALIZACION.id_institucion: fix description and PK/FK flags (FK, not PK).
$ws.Range("E45").Value = "Id de la institución que ofrece la especialización"
$ws.Range("F45").Value = "No"
$ws.Range("G45").Value = "Si"

# ESTUDIANTE_ESPECIALIZACION.estado: clarify description with an example.
$ws.Range("E63").Value = "Estado de la especialización para el estudiante. EJ: 'En curso'"

# ESTUDIANTE_CURSO.estado: clarify description with an example.
$ws.Range("E69").Value = "Estado actual del curso para el estudiante. Ej:'En curso'"

# ---------------------------------------------------------------------------
# 2) Row height tweaks
# ---------------------------------------------------------------------------

# Rows that go back to the default (non custom) height.
$ws.Rows.Item(36).RowHeight = 15
$ws.Rows.Item(36).UseStandardHeight = $true
$ws.Rows.Item(61).RowHeight = 15
$ws.Rows.Item(61).UseStandardHeight = $true
$ws.Rows.Item(62).RowHeight = 15
$ws.Rows.Item(62).UseStandardHeight = $true
$ws.Rows.Item(68).RowHeight = 15
$ws.Rows.Item(68).UseStandardHeight = $true
$ws.Rows.Item(77).RowHeight = 15
$ws.Rows.Item(77).UseStandardHeight = $true

# Rows that now need an explicit custom height.
$ws.Rows.Item(69).RowHeight = 38.25
$ws.Rows.Item(70).RowHeight = 42.75
$ws.Rows.Item(71).RowHeight = 25.5
$ws.Rows.Item(72).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 3) Column width tweak (column E got wider to fit the new text)
# ---------------------------------------------------------------------------

$ws.Columns.Item(5).ColumnWidth = 37.6

# ---------------------------------------------------------------------------
# 4) View state: zoom level and current selection
# ---------------------------------------------------------------------------

$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("C72").Select()
